# Enter results for "Match Amical 7" (columns BG=minutes, BH=T/R/HG status,
# BI=goals, BJ=assists) for every player row, and add the new player
# "Theo Owono" on row 30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ BG = minutes (optional); BH = status "T"/"R"/"HG"; BI = goals (optional); BJ = assists (optional) }
$matchData = @{
  2  = @{ BH = "HG" }
  3  = @{ BG = 45; BH = "R" }
  4  = @{ BG = 45; BH = "T" }
  5  = @{ BG = 45; BH = "R" }
  6  = @{ BH = "HG" }
  7  = @{ BG = 75; BH = "T" }
  8  = @{ BG = 45; BH = "R" }
  9  = @{ BH = "HG" }
  10 = @{ BH = "HG" }
  11 = @{ BG = 45; BH = "T" }
  13 = @{ BH = "HG" }
  14 = @{ BG = 90; BH = "T" }
  15 = @{ BG = 45; BH = "R" }
  16 = @{ BH = "HG" }
  17 = @{ BG = 45; BH = "T"; BI = 2 }
  18 = @{ BG = 45; BH = "T"; BJ = 2 }
  19 = @{ BG = 45; BH = "R" }
  20 = @{ BG = 45; BH = "T" }
  22 = @{ BG = 45; BH = "T" }
  24 = @{ BG = 60; BH = "T" }
  26 = @{ BG = 45; BH = "R" }
  27 = @{ BG = 45; BH = "T" }
  28 = @{ BG = 90; BH = "T" }
  29 = @{ BG = 45; BH = "R" }
}

foreach ($row in $matchData.Keys) {
  $info = $matchData[$row]
  if ($info.ContainsKey("BG")) {
    $ws.Cells.Item($row, 59).Value = $info["BG"]   # column BG = 59
  }
  if ($info.ContainsKey("BH")) {
    $ws.Cells.Item($row, 60).Value = $info["BH"]   # column BH = 60
  }
  if ($info.ContainsKey("BI")) {
    $ws.Cells.Item($row, 61).Value = $info["BI"]   # column BI = 61
  }
  if ($info.ContainsKey("BJ")) {
    $ws.Cells.Item($row, 62).Value = $info["BJ"]   # column BJ = 62
  }
}

# New player "Theo Owono" on row 30: name, the F-column "temps de jeu Matchs
# Amicaux" rollup formula (extended down from F29), and his Match Amical 7
# line (45 min, replacant).
$ws.Range("A30").Value = "Theo Owono"
$ws.Range("A30").HorizontalAlignment = -4108
$ws.Range("A30").VerticalAlignment = -4108

$ws.Range("F30").Formula = "=SUM(AI30,AM30,AQ30,AU30,AY30,BC30,BG30,BK30,BO30,BS30,BW30,CA30,CE30,CI30,CM30,CQ30,CU30)"

$ws.Range("BG30").Value = 45
$ws.Range("BH30").Value = "R"

# Extend the "R / B / P / RENFO-TEK / OK / NN" conditional formatting down
# through the newly added row.
$fc = $ws.Range("A22").FormatConditions
$fc.Item(1).ModifyAppliesToRange($ws.Range("A22:A30"))

# Reflect the author's final selection in the saved view.
$ws.Range("BQ22").Select()
